# feat: add 2022-Q1 data
#
# The workbook currently has 3 sheets: "2021-Q1", "2021-Q2", "总计".
# The "总计" sheet (sheet index 3) is repurposed into the new "2022-Q1"
# per-fund holdings sheet, and a duplicate of the original "总计" sheet is
# created right after it (preserving its formatting/styles exactly) and
# then updated with the new aggregate row for 2022-Q1.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item(2)   # "2021-Q2" - used as formatting reference
$wsTotal = $wb.Worksheets.Item(3)   # currently "总计"

# Duplicate the "总计" sheet so we keep an exact, fully-formatted copy that
# will become the new "总计" sheet, while the original object becomes
# "2022-Q1".
$wsTotal.Copy([System.Reflection.Missing]::Value, $wsTotal) | Out-Null

$wsQ1 = $wsTotal                       # will become "2022-Q1"
$wsNewTotal = $wb.Worksheets.Item(4)   # freshly created copy -> "总计"

$wsQ1.Name = "2022-Q1"
$wsNewTotal.Name = "总计"

# ---------------------------------------------------------------
# 1) Build the "2022-Q1" sheet (per-fund holdings), same layout as the
#    "2021-Q1" / "2021-Q2" sheets.
# ---------------------------------------------------------------

# Clear out the old "总计" style data rows (kept B2:D3 from old sheet) so we
# can lay out the new 7-column table cleanly.
$wsQ1.Range("B2:D3").ClearContents() | Out-Null

# Copy header formatting (bold font + border, style used on row 1 of the
# other per-fund sheets) across B1:H1.
$ws2.Range("B1:H1").Copy() | Out-Null
$wsQ1.Range("B1:H1").PasteSpecial(-4122) | Out-Null

# Copy the "A" column data-row style (centered bold w/ border) down to A3
# (A2 already carries it from the old sheet).
$ws2.Range("A2").Copy() | Out-Null
$wsQ1.Range("A3").PasteSpecial(-4122) | Out-Null

# Header row
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Row 2 (numeric-looking codes/figures must stay text, like the other
# sheets, so force text format, assign, then drop the now-unneeded
# "@" format so no stray style index is left on the cell).
$wsQ1.Range("A2").Value = 0

$wsQ1.Range("B2").NumberFormat = "@"
$wsQ1.Range("B2").Value = "398061"
$wsQ1.Range("B2").ClearFormats() | Out-Null

$wsQ1.Range("C2").Value = "中海消费混合"

$wsQ1.Range("D2").NumberFormat = "@"
$wsQ1.Range("D2").Value = "4.19"
$wsQ1.Range("D2").ClearFormats() | Out-Null

$wsQ1.Range("E2").NumberFormat = "@"
$wsQ1.Range("E2").Value = "93.41"
$wsQ1.Range("E2").ClearFormats() | Out-Null

$wsQ1.Range("F2").NumberFormat = "@"
$wsQ1.Range("F2").Value = "4.42"
$wsQ1.Range("F2").ClearFormats() | Out-Null

$wsQ1.Range("G2").NumberFormat = "@"
$wsQ1.Range("G2").Value = "0.1852"
$wsQ1.Range("G2").ClearFormats() | Out-Null

$wsQ1.Range("H2").Value = 2

# Row 3
$wsQ1.Range("A3").Value = 1

$wsQ1.Range("B3").NumberFormat = "@"
$wsQ1.Range("B3").Value = "180028"
$wsQ1.Range("B3").ClearFormats() | Out-Null

$wsQ1.Range("C3").Value = "银华永祥灵活配置混合"

$wsQ1.Range("D3").NumberFormat = "@"
$wsQ1.Range("D3").Value = "0.61"
$wsQ1.Range("D3").ClearFormats() | Out-Null

$wsQ1.Range("E3").NumberFormat = "@"
$wsQ1.Range("E3").Value = "77.23"
$wsQ1.Range("E3").ClearFormats() | Out-Null

$wsQ1.Range("F3").NumberFormat = "@"
$wsQ1.Range("F3").Value = "3.53"
$wsQ1.Range("F3").ClearFormats() | Out-Null

$wsQ1.Range("G3").NumberFormat = "@"
$wsQ1.Range("G3").Value = "0.0215"
$wsQ1.Range("G3").ClearFormats() | Out-Null

$wsQ1.Range("H3").Value = 6

# ---------------------------------------------------------------
# 2) Build the new "总计" sheet: same as before, plus a new first data row
#    for 2022-Q1, with older rows shifted down.
# ---------------------------------------------------------------

# Give row 4 the same "A" column style as the other data rows (A2/A3
# already have it from the duplicated sheet).
$wsNewTotal.Range("A2").Copy() | Out-Null
$wsNewTotal.Range("A4").PasteSpecial(-4122) | Out-Null

$wsNewTotal.Range("A4").Value = 2
$wsNewTotal.Range("B4").Value = "2021-Q1"
$wsNewTotal.Range("C4").Value = 4
$wsNewTotal.Range("D4").Value = 0.26

$wsNewTotal.Range("A3").Value = 1
$wsNewTotal.Range("B3").Value = "2021-Q2"
$wsNewTotal.Range("C3").Value = 6
$wsNewTotal.Range("D3").Value = 0.59

$wsNewTotal.Range("A2").Value = 0
$wsNewTotal.Range("B2").Value = "2022-Q1"
$wsNewTotal.Range("C2").Value = 2
$wsNewTotal.Range("D2").Value = 0.21

# Restore the originally active/selected sheet so we don't leave the
# duplicated sheet marked as the active tab.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate() | Out-Null
